$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) EvaJourney sheet (sheet1) — wording tweaks on a few PPV-delay lines
# ---------------------------------------------------------------------------
$eva = $wb.Worksheets.Item("EvaJourney")
$eva.Range("B4").Value  = "I'm about to finish... stay with me"
$eva.Range("B5").Value  = "don't go anywhere"
$eva.Range("B11").Value = "screw it I'm done waiting"

# ---------------------------------------------------------------------------
# 2) Rename "cumcontrol" -> "cumcontrol1" and refresh its copy text
# ---------------------------------------------------------------------------
$cc1 = $wb.Worksheets.Item("cumcontrol")
$cc1.Name = "cumcontrol1"

$cc1.Range("B2").Value = "just hold on a little more, I want the last thing you see to be this"

$cc1.Range("B3").Value = "wait for me papi... I have one more thing and I want you to see it before we finish"
$cc1.Range("C3").Value = "DELAY. Send PPV."

$cc1.Range("B4").Value = "stay with me, I'm almost there too... watch this"
$cc1.Range("C4").Value = "SYNC variant. Send PPV."

$cc1.Range("B5").Value = "I want us to finish together papi... open this and let go with me"
$cc1.Range("C5").Value = "SYNC. Send PPV."

$cc1.Range("B6").Value = "please don't finish yet... I'm not ready for this to be over"

$cc1.Range("B7").Value = "not yet papi... I want this to last a little longer with you"
$cc1.Range("C7").Value = "CONTROL."

# ---------------------------------------------------------------------------
# 3) Insert a brand new "cumcontrol2" sheet right after "cumcontrol1" (and
#    before "dickpic"), cloning the cumcontrol1 layout/formatting, then
#    overwrite it with its own fresh copy set.
# ---------------------------------------------------------------------------
$cc2 = $wb.Worksheets.Add($null, $cc1)
$cc2.Name = "cumcontrol2"

$cc1.Range("A1:D7").Copy($cc2.Range("A1"))

$cc2.Columns.Item(1).ColumnWidth = 19.166666666666668
$cc2.Columns.Item(2).ColumnWidth = 79.16666666666667
$cc2.Columns.Item(3).ColumnWidth = 49.166666666666664
$cc2.Columns.Item(4).ColumnWidth = 24.166666666666668

$cc2.Range("A2").Value = "delay2"
$cc2.Range("B2").Value = "just a little longer for me papi? the next one is special"
$cc2.Range("C2").Value = "DELAY variant."

$cc2.Range("A3").Value = "delay1"
$cc2.Range("B3").Value = "please wait... what I'm about to send, I want you to really take it in"
$cc2.Range("C3").Value = "DELAY. Send PPV."

$cc2.Range("A4").Value = "sync2"
$cc2.Range("B4").Value = "I need you to see this before we both let go"
$cc2.Range("C4").Value = "SYNC variant."

$cc2.Range("A5").Value = "sync1"
$cc2.Range("B5").Value = "okay papi... together, right now... open this"
$cc2.Range("C5").Value = "SYNC. Send PPV."

$cc2.Range("A6").Value = "edge2"
$cc2.Range("B6").Value = "don't rush... this is too good to end yet"
$cc2.Range("C6").Value = "EDGE variant."

$cc2.Range("A7").Value = "edge1"
$cc2.Range("B7").Value = "slow down papi... I want to feel every second of this with you"
$cc2.Range("C7").Value = "CONTROL."

Write-Output "done"
